$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 140 - this shifts the existing rows 140:217
# down to 141:218, preserving all their values (incl. the constant A/B/C/
# E/F/G/H/I/K/L/M/N/O/P/Q/R pattern) automatically.
$ws.Rows("140:140").Insert()

# Populate the newly inserted (blank) row 140 with its record. All columns
# other than D (Fecha) and J (Volumen) follow the same constant pattern as
# every other row in this block.
$ws.Range("A140").Value = 3
$ws.Range("B140").Value = "Femacal de La Calera"
$ws.Range("C140").Value = "Coquimbo"
$ws.Range("D140").Value = 44529
$ws.Range("E140").Value = 5
$ws.Range("F140").Value = 100112039
$ws.Range("G140").Value = "Ciboulette"
$ws.Range("H140").Value = "Sin especificar"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 160
$ws.Range("K140").Value = 1500
$ws.Range("L140").Value = 1500
$ws.Range("M140").Value = 1500
$ws.Range("N140").Value = '$/docena de atados'
$ws.Range("O140").Value = "Provincia de Quillota"
$ws.Range("P140").Value = 500
$ws.Range("Q140").Value = 3
$ws.Range("R140").Value = "Hortaliza"
